# Add two new columns, I ("I0") and J ("IF"), to Sheet1, with header labels
# in row 1 and numeric data in rows 2-40 (mirrors the existing H column's
# header style/border).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
# Copy H1's formatting (bold font + border + centered alignment) onto I1:J1
# before setting their text, so the new header cells match the look of the
# existing B1:H1 headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data (rows 2-40) ---
# row number -> (I value, J value)
$data = @{
    2  = @(2, 3)
    3  = @(4, 6)
    4  = @(6, 6)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(1, 3)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 5)
    11 = @(1, 5)
    12 = @(1, 5)
    13 = @(1, 3)
    14 = @(6, 7)
    15 = @(1, 7)
    16 = @(1, 6)
    17 = @(1, 7)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 4)
    21 = @(1, 4)
    22 = @(1, 6)
    23 = @(1, 6)
    24 = @(1, 7)
    25 = @(1, 6)
    26 = @(1, 5)
    27 = @(1, 7)
    28 = @(1, 6)
    29 = @(1, 6)
    30 = @(1, 5)
    31 = @(1, 5)
    32 = @(1, 7)
    33 = @(1, 5)
    34 = @(1, 7)
    35 = @(1, 7)
    36 = @(1, 5)
    37 = @(1, 4)
    38 = @(1, 3)
    39 = @(1, 2)
    40 = @(1, 1)
}

foreach ($r in $data.Keys) {
    $pair = $data[$r]
    $ws.Cells.Item($r, 9).Value  = $pair[0]   # column I
    $ws.Cells.Item($r, 10).Value = $pair[1]   # column J
}
